# Re-apply the betexplorer scraper refresh for
# 2023/indonesia_liga-1_2023-2024.xlsx
#
# 1) A handful of same-day doubleheader rows got re-ordered by the
#    scraper (the two matches on a given date swapped positions); the
#    visible effect is that columns F:V (home team .. url) of each pair
#    of adjacent rows are exchanged, while A:E (index/country/tourney/
#    season/date) stay put.
# 2) Nine new match rows (145-153) were appended at the bottom, for
#    matches played 25-30/10/2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Part 1: swap columns F..V between each of these adjacent row pairs
# ---------------------------------------------------------------
$swapCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
$swapPairs = @(
    @(28,29),
    @(30,31),
    @(37,38),
    @(42,43),
    @(47,48),
    @(51,52),
    @(57,58),
    @(60,61),
    @(62,63),
    @(76,77),
    @(78,79),
    @(82,83),
    @(98,99),
    @(134,135)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $swapCols) {
        $ref1 = "$col$r1"
        $ref2 = "$col$r2"
        $v1 = $ws.Range($ref1).Value2
        $v2 = $ws.Range($ref2).Value2
        $ws.Range($ref1).Value = $v2
        $ws.Range($ref2).Value = $v1
    }
}

# ---------------------------------------------------------------
# Part 2: append the nine new rows (145..153), with A=index (s=1 style)
#         and E=match-datetime (s=2 style) copied from the formats used
#         by the existing rows.
# ---------------------------------------------------------------
$newRows = @(
    @{ A=144; E=45226.41666666666; F="Persik Kediri";  G=4; H="Persebaya";       I=0; J=2.1;  K="25/10/2023 22:12"; L=2.61; M="27/10/2023 09:58"; N=3.35; O="25/10/2023 22:12"; P=3.55; Q="27/10/2023 09:52"; R=3.06; S="25/10/2023 22:12"; T=2.54; U="27/10/2023 09:56"; V="https://www.betexplorer.com/football/indonesia/liga-1/persik-kediri-persebaya/CGkUL9KD/" },
    @{ A=145; E=45226.58333333334; F="Barito Putera";  G=1; H="Persikabo 1973";  I=1; J=1.45; K="26/10/2023 02:12"; L=1.35; M="27/10/2023 13:58"; N=4.29; O="26/10/2023 02:12"; P=4.87; Q="27/10/2023 13:58"; R=5.51; S="26/10/2023 02:12"; T=8.98; U="27/10/2023 13:58"; V="https://www.betexplorer.com/football/indonesia/liga-1/ps-barito-putera-persikabo-1973/YPjQMk57/" },
    @{ A=146; E=45227.41666666666; F="Arema FC";       G=1; H="Madura United";   I=1; J=3.06; K="26/10/2023 22:12"; L=3.81; M="28/10/2023 09:57"; N=3.2;  O="26/10/2023 22:12"; P=3.44; Q="28/10/2023 09:57"; R=2.17; S="26/10/2023 22:12"; T=1.98; U="28/10/2023 09:57"; V="https://www.betexplorer.com/football/indonesia/liga-1/arema-fc-madura-united/AiRMMXKD/" },
    @{ A=147; E=45227.58333333334; F="Borneo";         G=3; H="Dewa United";     I=1; J=1.74; K="27/10/2023 02:13"; L=1.59; M="28/10/2023 13:21"; N=3.61; O="27/10/2023 02:13"; P=3.99; Q="28/10/2023 13:21"; R=4.2;  S="27/10/2023 02:13"; T=5.64; U="28/10/2023 13:21"; V="https://www.betexplorer.com/football/indonesia/liga-1/borneo-dewa-united/KQGdGVcs/" },
    @{ A=148; E=45227.58333333334; F="Persib Bandung"; G=4; H="PSS Sleman";      I=1; J=1.35; K="27/10/2023 02:13"; L=1.38; M="28/10/2023 13:58"; N=4.75; O="27/10/2023 02:13"; P=4.66; Q="28/10/2023 13:58"; R=6.7;  S="27/10/2023 02:13"; T=8.38; U="28/10/2023 13:58"; V="https://www.betexplorer.com/football/indonesia/liga-1/persib-bandung-pss-sleman/OfVQLizK/" },
    @{ A=149; E=45228.375;         F="Persis Solo";    G=2; H="FC Bhayangkara";  I=1; J=1.72; K="27/10/2023 22:13"; L=1.53; M="29/10/2023 08:55"; N=3.74; O="27/10/2023 22:13"; P=4.36; Q="29/10/2023 08:58"; R=3.92; S="27/10/2023 22:13"; T=5.71; U="29/10/2023 08:58"; V="https://www.betexplorer.com/football/indonesia/liga-1/persis-solo-fc-bhayangkara/IHUUKBkQ/" },
    @{ A=150; E=45228.54166666666; F="PSIS Semarang";  G=2; H="Persija Jakarta"; I=1; J=2.49; K="28/10/2023 02:13"; L=1.83; M="29/10/2023 12:54"; N=3.12; O="28/10/2023 02:13"; P=3.38; Q="29/10/2023 12:54"; R=2.71; S="28/10/2023 02:13"; T=4.69; U="29/10/2023 12:54"; V="https://www.betexplorer.com/football/indonesia/liga-1/psis-semarang-persija-jakarta/x2B1FkDm/" },
    @{ A=151; E=45229.375;         F="RANS Nusantara"; G=1; H="PSM Makassar";    I=1; J=2.73; K="28/10/2023 22:13"; L=2.37; M="30/10/2023 08:57"; N=2.8;  O="28/10/2023 22:13"; P=3.22; Q="30/10/2023 08:59"; R=2.71; S="28/10/2023 22:13"; T=3.08; U="30/10/2023 08:57"; V="https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-psm-makassar/GhA5E9Sg/" },
    @{ A=152; E=45229.54166666666; F="Bali United";    G=3; H="Persita";         I=0; J=1.55; K="29/10/2023 02:13"; L=1.47; M="30/10/2023 12:23"; N=4.03; O="29/10/2023 02:13"; P=4.06; Q="30/10/2023 12:58"; R=4.83; S="29/10/2023 02:13"; T=6;    U="30/10/2023 12:58"; V="https://www.betexplorer.com/football/indonesia/liga-1/bali-united-persita/AF99DTs0/" }
)

$startRow = 145
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    # Pull formatting (bold/border/centered for A, datetime numfmt for E)
    # from the last existing data row so the new cells get the same
    # style slot instead of a brand-new one.
    $ws.Cells.Item(144, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item(144, 5).Copy()
    $ws.Cells.Item($row, 5).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $data.A
    $ws.Cells.Item($row, 2).Value = "indonesia"
    $ws.Cells.Item($row, 3).Value = "liga-1"
    $ws.Cells.Item($row, 4).Value = "2023-2024"
    $ws.Cells.Item($row, 5).Value = $data.E
    $ws.Cells.Item($row, 6).Value = $data.F
    $ws.Cells.Item($row, 7).Value = $data.G
    $ws.Cells.Item($row, 8).Value = $data.H
    $ws.Cells.Item($row, 9).Value = $data.I
    $ws.Cells.Item($row, 10).Value = $data.J
    $ws.Cells.Item($row, 11).Value = $data.K
    $ws.Cells.Item($row, 12).Value = $data.L
    $ws.Cells.Item($row, 13).Value = $data.M
    $ws.Cells.Item($row, 14).Value = $data.N
    $ws.Cells.Item($row, 15).Value = $data.O
    $ws.Cells.Item($row, 16).Value = $data.P
    $ws.Cells.Item($row, 17).Value = $data.Q
    $ws.Cells.Item($row, 18).Value = $data.R
    $ws.Cells.Item($row, 19).Value = $data.S
    $ws.Cells.Item($row, 20).Value = $data.T
    $ws.Cells.Item($row, 21).Value = $data.U
    $ws.Cells.Item($row, 22).Value = $data.V
}
